$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-16 Saturday" "2024-03-17 Sunday"
Replace-Text "745÷4=186, 1" "170÷9=18, 8"
Replace-Text "786÷6=131, 0" "981÷9=109, 0"
Replace-Text "979÷4=244, 3" "753÷2=376, 1"
Replace-Text "594÷4=148, 2" "948÷5=189, 3"
Replace-Text "805÷4=201, 1" "820÷2=410, 0"
Replace-Text "848÷2=424, 0" "374÷5=74, 4"
Replace-Text "273÷3=91, 0" "798÷7=114, 0"
Replace-Text "173÷9=19, 2" "627÷4=156, 3"
Replace-Text "540÷7=77, 1" "979÷8=122, 3"
Replace-Text "783÷7=111, 6" "322÷8=40, 2"
Replace-Text "361÷8=45, 1" "755÷3=251, 2"
Replace-Text "952÷3=317, 1" "328÷6=54, 4"
Replace-Text "506÷6=84, 2" "253÷7=36, 1"
Replace-Text "927÷4=231, 3" "835÷3=278, 1"
Replace-Text "527÷7=75, 2" "854÷7=122, 0"
Replace-Text "892÷5=178, 2" "604÷4=151, 0"
Replace-Text "799÷7=114, 1" "917÷8=114, 5"
Replace-Text "956÷2=478, 0" "710÷7=101, 3"
Replace-Text "241÷5=48, 1" "382÷5=76, 2"
Replace-Text "235÷2=117, 1" "790÷6=131, 4"
Replace-Text "383÷6=63, 5" "252÷6=42, 0"
Replace-Text "131÷7=18, 5" "431÷7=61, 4"
Replace-Text "165÷7=23, 4" "908÷8=113, 4"
Replace-Text "104÷6=17, 2" "492÷3=164, 0"
Replace-Text "809÷5=161, 4" "422÷7=60, 2"

Write-Host "Done"
